$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("meta")

# Insert a new row before the existing "style"/"default" row (row 9), shifting it
# (and the trailing blank row) down by one.
$ws.Rows.Item(9).Insert()

$ws.Range("A9").Value = "x_date_format"

# Use a formula + paste-as-values round trip so Excel stores the literal leading
# apostrophe as real text ("'yy") instead of treating it as a quote-prefix marker
# (which would silently drop the apostrophe from the stored string).
$ws.Range("B9").Formula = "=""'yy"""
$ws.Range("B9").Copy()
$ws.Range("B9").PasteSpecial(-4163)  # xlPasteValues

# Copy the "key" column formatting (bold/orange) from the row above onto the new cell.
$ws.Range("A8").Copy()
$ws.Range("A9").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = 0
